# Auto-generated Excel COM-interop script
# Applies cell-level numeric updates to the Leve profit tracking sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 778
$ws.Cells.Item(15, 9).Value = 778
$ws.Cells.Item(15, 11).Value = 2334
$ws.Cells.Item(15, 13).Value = -2165
$ws.Cells.Item(28, 8).Value = 3874.4285
$ws.Cells.Item(28, 9).Value = 946.5
$ws.Cells.Item(28, 11).Value = 946.5
$ws.Cells.Item(28, 13).Value = -461.5
$ws.Cells.Item(29, 8).Value = 3975.1365
$ws.Cells.Item(29, 10).Value = 7490
$ws.Cells.Item(29, 12).Value = 22470
$ws.Cells.Item(29, 14).Value = -23032
$ws.Cells.Item(74, 8).Value = 2923.5
$ws.Cells.Item(74, 10).Value = 1199
$ws.Cells.Item(74, 12).Value = 1199
$ws.Cells.Item(74, 14).Value = -3071
$ws.Cells.Item(77, 8).Value = 2923.5
$ws.Cells.Item(77, 10).Value = 1199
$ws.Cells.Item(77, 12).Value = 5995
$ws.Cells.Item(77, 14).Value = -15355
$ws.Cells.Item(92, 8).Value = 139.44444
$ws.Cells.Item(92, 9).Value = 158
$ws.Cells.Item(92, 11).Value = 158
$ws.Cells.Item(92, 13).Value = 1090
$ws.Cells.Item(111, 8).Value = 955.7143
$ws.Cells.Item(111, 10).Value = 1733.3334
$ws.Cells.Item(111, 12).Value = 5200.0002
$ws.Cells.Item(111, 14).Value = -11334.0002
$ws.Cells.Item(116, 8).Value = 3500
$ws.Cells.Item(116, 9).Value = 3000
$ws.Cells.Item(116, 10).Value = 4000
$ws.Cells.Item(116, 11).Value = 3000
$ws.Cells.Item(116, 12).Value = 4000
$ws.Cells.Item(116, 13).Value = 442
$ws.Cells.Item(116, 14).Value = -10884
$ws.Cells.Item(137, 8).Value = 1964.1666
$ws.Cells.Item(137, 9).Value = 1357
$ws.Cells.Item(137, 11).Value = 4071
$ws.Cells.Item(137, 13).Value = -1521
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2445.1
$ws.Cells.Item(45, 9).Value = 1594.7142
$ws.Cells.Item(45, 10).Value = 4429.3335
$ws.Cells.Item(45, 11).Value = 1594.7142
$ws.Cells.Item(45, 12).Value = 4429.3335
$ws.Cells.Item(45, 13).Value = -1217.7142
$ws.Cells.Item(45, 14).Value = -5183.3335
$ws.Cells.Item(61, 8).Value = 8250
$ws.Cells.Item(61, 9).Value = 10000
$ws.Cells.Item(61, 10).Value = 8000
$ws.Cells.Item(61, 11).Value = 10000
$ws.Cells.Item(61, 12).Value = 8000
$ws.Cells.Item(61, 13).Value = -9788
$ws.Cells.Item(61, 14).Value = -8424
$ws.Cells.Item(63, 8).Value = 9389.799999999999
$ws.Cells.Item(63, 9).Value = 1818.8
$ws.Cells.Item(63, 10).Value = 16960.8
$ws.Cells.Item(63, 11).Value = 1818.8
$ws.Cells.Item(63, 12).Value = 16960.8
$ws.Cells.Item(63, 13).Value = -1132.8
$ws.Cells.Item(63, 14).Value = -18332.8
$ws.Cells.Item(66, 8).Value = 9389.799999999999
$ws.Cells.Item(66, 9).Value = 1818.8
$ws.Cells.Item(66, 10).Value = 16960.8
$ws.Cells.Item(66, 11).Value = 9094
$ws.Cells.Item(66, 12).Value = 84804
$ws.Cells.Item(66, 13).Value = -5662
$ws.Cells.Item(66, 14).Value = -91668
$ws.Cells.Item(74, 8).Value = 2472.5334
$ws.Cells.Item(74, 9).Value = 2472.5334
$ws.Cells.Item(74, 11).Value = 2472.5334
$ws.Cells.Item(74, 13).Value = -1598.5334
$ws.Cells.Item(77, 8).Value = 2472.5334
$ws.Cells.Item(77, 9).Value = 2472.5334
$ws.Cells.Item(77, 11).Value = 12362.667
$ws.Cells.Item(77, 13).Value = -7994.666999999999
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 13).ClearContents()
$ws.Cells.Item(132, 8).Value = 5620.2856
$ws.Cells.Item(132, 9).Value = 6052.2
$ws.Cells.Item(132, 11).Value = 18156.6
$ws.Cells.Item(132, 13).Value = -15626.6
$ws.Cells.Item(136, 8).Value = 8250
$ws.Cells.Item(136, 9).Value = 10000
$ws.Cells.Item(136, 10).Value = 8000
$ws.Cells.Item(136, 11).Value = 30000
$ws.Cells.Item(136, 12).Value = 24000
$ws.Cells.Item(136, 13).Value = -27450
$ws.Cells.Item(136, 14).Value = -29100
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 969.9091
$ws.Cells.Item(3, 9).Value = 863.2222
$ws.Cells.Item(3, 11).Value = 863.2222
$ws.Cells.Item(3, 13).Value = -749.2222
$ws.Cells.Item(86, 8).Value = 5414.409
$ws.Cells.Item(86, 9).Value = 4271.4614
$ws.Cells.Item(86, 10).Value = 7065.3335
$ws.Cells.Item(86, 11).Value = 4271.4614
$ws.Cells.Item(86, 12).Value = 7065.3335
$ws.Cells.Item(86, 13).Value = -3148.4614
$ws.Cells.Item(86, 14).Value = -9311.333500000001
$ws.Cells.Item(89, 8).Value = 5414.409
$ws.Cells.Item(89, 9).Value = 4271.4614
$ws.Cells.Item(89, 10).Value = 7065.3335
$ws.Cells.Item(89, 11).Value = 21357.307
$ws.Cells.Item(89, 12).Value = 35326.6675
$ws.Cells.Item(89, 13).Value = -15741.307
$ws.Cells.Item(89, 14).Value = -46558.6675
$ws.Cells.Item(134, 8).Value = 3337.375
$ws.Cells.Item(134, 9).Value = 3337.375
$ws.Cells.Item(134, 11).Value = 10012.125
$ws.Cells.Item(134, 13).Value = -7477.125
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 460.6111
$ws.Cells.Item(7, 9).Value = 110.07143
$ws.Cells.Item(7, 10).Value = 1687.5
$ws.Cells.Item(7, 11).Value = 110.07143
$ws.Cells.Item(7, 12).Value = 1687.5
$ws.Cells.Item(7, 13).Value = 2.928569999999993
$ws.Cells.Item(7, 14).Value = -1913.5
$ws.Cells.Item(31, 8).Value = 5776.1934
$ws.Cells.Item(31, 9).Value = 2630.2307
$ws.Cells.Item(31, 11).Value = 2630.2307
$ws.Cells.Item(31, 13).Value = -2335.2307
$ws.Cells.Item(34, 8).Value = 5776.1934
$ws.Cells.Item(34, 9).Value = 2630.2307
$ws.Cells.Item(34, 11).Value = 2630.2307
$ws.Cells.Item(34, 13).Value = -2428.2307
$ws.Cells.Item(99, 8).Value = 2343.9092
$ws.Cells.Item(99, 9).Value = 1988.6666
$ws.Cells.Item(99, 11).Value = 1988.6666
$ws.Cells.Item(99, 13).Value = -490.6666
$ws.Cells.Item(126, 8).Value = 2343.9092
$ws.Cells.Item(126, 9).Value = 1988.6666
$ws.Cells.Item(126, 11).Value = 5965.9998
$ws.Cells.Item(126, 13).Value = -3495.9998
$ws.Cells.Item(132, 8).Value = 996
$ws.Cells.Item(132, 9).Value = 996
$ws.Cells.Item(132, 11).Value = 2988
$ws.Cells.Item(132, 13).Value = -458
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 2179
$ws.Cells.Item(132, 9).Value = 1949.5
$ws.Cells.Item(132, 10).Value = 2199.8635
$ws.Cells.Item(132, 11).Value = 17545.5
$ws.Cells.Item(132, 12).Value = 19798.7715
$ws.Cells.Item(132, 13).Value = -15015.5
$ws.Cells.Item(132, 14).Value = -24858.7715
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 464.16666
$ws.Cells.Item(97, 9).Value = 517
$ws.Cells.Item(97, 10).Value = 200
$ws.Cells.Item(97, 11).Value = 517
$ws.Cells.Item(97, 12).Value = 200
$ws.Cells.Item(97, 13).Value = -21
$ws.Cells.Item(97, 14).Value = -1192
$ws.Cells.Item(102, 8).Value = 1579.0625
$ws.Cells.Item(102, 9).Value = 1351.0667
$ws.Cells.Item(102, 11).Value = 1351.0667
$ws.Cells.Item(102, 13).Value = 270.9332999999999
$ws.Cells.Item(113, 8).Value = 8727.299999999999
$ws.Cells.Item(113, 9).Value = 7879.1665
$ws.Cells.Item(113, 11).Value = 7879.1665
$ws.Cells.Item(113, 13).Value = -5709.1665
$ws.Cells.Item(122, 8).Value = 3995
$ws.Cells.Item(122, 9).Value = 2104
$ws.Cells.Item(122, 11).Value = 6312
$ws.Cells.Item(122, 13).Value = -3862
$ws.Cells.Item(132, 8).Value = 10000
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 10000
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 30000
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).Value = -35060
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 3020.5
$ws.Cells.Item(61, 9).Value = 886.4286
$ws.Cells.Item(61, 11).Value = 886.4286
$ws.Cells.Item(61, 13).Value = -684.4286
$ws.Cells.Item(113, 8).Value = 3020.5
$ws.Cells.Item(113, 9).Value = 886.4286
$ws.Cells.Item(113, 11).Value = 886.4286
$ws.Cells.Item(113, 13).Value = 1283.5714
$ws.Cells.Item(122, 8).Value = 985
$ws.Cells.Item(122, 9).Value = 985
$ws.Cells.Item(122, 11).Value = 2955
$ws.Cells.Item(122, 13).Value = -505
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 13).ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1914.7222
$ws.Cells.Item(132, 9).Value = 1850.8823
$ws.Cells.Item(132, 10).Value = 3000
$ws.Cells.Item(132, 11).Value = 5552.6469
$ws.Cells.Item(132, 12).Value = 9000
$ws.Cells.Item(132, 13).Value = -3022.6469
$ws.Cells.Item(132, 14).Value = -14060

Write-Output "Applied all profit/price updates"